$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Num($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Set-Txt($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
}

function Set-Bool($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Clear-Cell($row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

# ---------------------------------------------------------------
# Row 10 / Row 12 value swap (A, P, Q, R, AC columns)
# ---------------------------------------------------------------

# Row 10 becomes the old row-12 location data
Set-Num 10 1  111528300
Set-Txt 10 16 "Fläcksberget, Hjd"
Set-Num 10 17 467795.2212022893
Set-Num 10 18 6875452.272210476
Set-Txt 10 29 "Tre blommande."

# Row 12 becomes the old row-10 location data
Set-Num 12 1  111527876
Set-Txt 12 16 "Fläcksberget V, Hjd"
Set-Num 12 17 467615.2905344999
Set-Num 12 18 6875426.740629551
Clear-Cell 12 29

# ---------------------------------------------------------------
# Row 14 replaced with a new species record (Blå taggsvamp)
# ---------------------------------------------------------------

Set-Num 14 1  111908364
Set-Num 14 2  90660
Set-Txt 14 3  "Ovaliderad"
Set-Txt 14 4  "NT"
Set-Num 14 5  4362
Set-Txt 14 6  "Blå taggsvamp"
Set-Txt 14 7  "Hydnellum caeruleum"
Set-Txt 14 8  "(Hornem.) P.Karst."
Set-Txt 14 16 "Gröbäcken, Hjd"
Set-Num 14 17 467724.2196293612
Set-Num 14 18 6874811.291555981
Set-Num 14 19 20
Set-Txt 14 20 "Jämtland"
Set-Txt 14 21 "Härjedalen"
Set-Txt 14 22 "Härjedalen"
Set-Txt 14 23 "Sveg"
Set-Txt 14 25 "2023-09-05"
Set-Txt 14 26 "00:00"
Set-Txt 14 27 "2023-09-05"
Set-Txt 14 28 "00:00"
Set-Bool 14 30 $false
Set-Bool 14 31 $false
Set-Bool 14 33 $false
Set-Txt 14 49 "lennart karlsson"
Set-Txt 14 50 "lennart karlsson"

# ---------------------------------------------------------------
# New row 15 : Kolflarnlav
# ---------------------------------------------------------------

Set-Num 15 1  111909536
Set-Num 15 2  77267
Set-Txt 15 3  "Ovaliderad"
Set-Txt 15 4  "NT"
Set-Num 15 5  6446
Set-Txt 15 6  "Kolflarnlav"
Set-Txt 15 7  "Carbonicola anthracophila"
Set-Txt 15 8  "(Nyl.) Bendiksby & Timdal"
Set-Txt 15 16 "Fläcksberget, Hjd"
Set-Num 15 17 467891.3929605001
Set-Num 15 18 6875425.059267788
Set-Num 15 19 20
Set-Txt 15 20 "Jämtland"
Set-Txt 15 21 "Härjedalen"
Set-Txt 15 22 "Härjedalen"
Set-Txt 15 23 "Sveg"
Set-Txt 15 25 "2023-09-05"
Set-Txt 15 26 "00:00"
Set-Txt 15 27 "2023-09-05"
Set-Txt 15 28 "00:00"
Set-Bool 15 30 $false
Set-Bool 15 31 $false
Set-Bool 15 33 $false
Set-Txt 15 49 "lennart karlsson"
Set-Txt 15 50 "lennart karlsson"

# ---------------------------------------------------------------
# New row 16 : Doftticka
# ---------------------------------------------------------------

Set-Num 16 1  111908700
Set-Num 16 2  89965
Set-Txt 16 3  "Ovaliderad"
Set-Txt 16 4  "VU"
Set-Num 16 5  760
Set-Txt 16 6  "Doftticka"
Set-Txt 16 7  "Haploporus odorus"
Set-Txt 16 8  "(Sommerf.) Bondartsev & Singer"
Set-Txt 16 9  "6"
Set-Txt 16 10 "fruktkroppar"
Set-Txt 16 16 "Fläcksberget, Hjd"
Set-Num 16 17 467921.7931363151
Set-Num 16 18 6875306.87748003
Set-Num 16 19 20
Set-Txt 16 20 "Jämtland"
Set-Txt 16 21 "Härjedalen"
Set-Txt 16 22 "Härjedalen"
Set-Txt 16 23 "Sveg"
Set-Txt 16 25 "2023-09-05"
Set-Txt 16 26 "00:00"
Set-Txt 16 27 "2023-09-05"
Set-Txt 16 28 "00:00"
Set-Txt 16 29 "Förekomst av doftticka i avverkningsanmält område."
Set-Bool 16 30 $false
Set-Bool 16 31 $false
Set-Bool 16 33 $false
Set-Txt 16 49 "lennart karlsson"
Set-Txt 16 50 "lennart karlsson"

# ---------------------------------------------------------------
# New row 17 : Knärot
# ---------------------------------------------------------------

Set-Num 17 1  111908768
Set-Num 17 2  96348
Set-Txt 17 3  "Ovaliderad"
Set-Txt 17 4  "VU"
Set-Num 17 5  220787
Set-Txt 17 6  "Knärot"
Set-Txt 17 7  "Goodyera repens"
Set-Txt 17 8  "(L.) R. Br."
Set-Txt 17 9  "1"
Set-Txt 17 16 "Fläcksberget, Hjd"
Set-Num 17 17 467911.8445363804
Set-Num 17 18 6875299.456096188
Set-Num 17 19 20
Set-Txt 17 20 "Jämtland"
Set-Txt 17 21 "Härjedalen"
Set-Txt 17 22 "Härjedalen"
Set-Txt 17 23 "Sveg"
Set-Txt 17 25 "2023-09-05"
Set-Txt 17 26 "00:00"
Set-Txt 17 27 "2023-09-05"
Set-Txt 17 28 "00:00"
Set-Bool 17 30 $false
Set-Bool 17 31 $false
Set-Bool 17 33 $false
Set-Txt 17 49 "lennart karlsson"
Set-Txt 17 50 "lennart karlsson"

# ---------------------------------------------------------------
# New row 18 : Rödgul trumpetsvamp
# ---------------------------------------------------------------

Set-Num 18 1  111909766
Set-Num 18 2  89183
Set-Txt 18 3  "Ovaliderad"
Set-Txt 18 4  "LC"
Set-Num 18 5  3215
Set-Txt 18 6  "Rödgul trumpetsvamp"
Set-Txt 18 7  "Craterellus lutescens"
Set-Txt 18 8  "(Fr.) Fr."
Set-Txt 18 16 "Fläcksberget, Hjd"
Set-Num 18 17 467756.8135427741
Set-Num 18 18 6875469.545251801
Set-Num 18 19 20
Set-Txt 18 20 "Jämtland"
Set-Txt 18 21 "Härjedalen"
Set-Txt 18 22 "Härjedalen"
Set-Txt 18 23 "Sveg"
Set-Txt 18 25 "2023-09-05"
Set-Txt 18 26 "00:00"
Set-Txt 18 27 "2023-09-05"
Set-Txt 18 28 "00:00"
Set-Bool 18 30 $false
Set-Bool 18 31 $false
Set-Bool 18 33 $false
Set-Txt 18 49 "lennart karlsson"
Set-Txt 18 50 "lennart karlsson"

# ---------------------------------------------------------------
# New row 19 : Kolflarnlav
# ---------------------------------------------------------------

Set-Num 19 1  111909174
Set-Num 19 2  77267
Set-Txt 19 3  "Ovaliderad"
Set-Txt 19 4  "NT"
Set-Num 19 5  6446
Set-Txt 19 6  "Kolflarnlav"
Set-Txt 19 7  "Carbonicola anthracophila"
Set-Txt 19 8  "(Nyl.) Bendiksby & Timdal"
Set-Txt 19 16 "Fläcksberget, Hjd"
Set-Num 19 17 467989.0228066717
Set-Num 19 18 6875352.744105402
Set-Num 19 19 20
Set-Txt 19 20 "Jämtland"
Set-Txt 19 21 "Härjedalen"
Set-Txt 19 22 "Härjedalen"
Set-Txt 19 23 "Sveg"
Set-Txt 19 25 "2023-09-05"
Set-Txt 19 26 "00:00"
Set-Txt 19 27 "2023-09-05"
Set-Txt 19 28 "00:00"
Set-Bool 19 30 $false
Set-Bool 19 31 $false
Set-Bool 19 33 $false
Set-Txt 19 49 "lennart karlsson"
Set-Txt 19 50 "lennart karlsson"

# ---------------------------------------------------------------
# New row 20 : Kolflarnlav (this was the original row 14 data)
# ---------------------------------------------------------------

Set-Num 20 1  111942712
Set-Num 20 2  77267
Set-Txt 20 3  "Ovaliderad"
Set-Txt 20 4  "NT"
Set-Num 20 5  6446
Set-Txt 20 6  "Kolflarnlav"
Set-Txt 20 7  "Carbonicola anthracophila"
Set-Txt 20 8  "(Nyl.) Bendiksby & Timdal"
Set-Txt 20 16 "Fläcksberget, Hjd"
Set-Num 20 17 468231.4750461024
Set-Num 20 18 6875021.661872049
Set-Num 20 19 20
Set-Txt 20 20 "Jämtland"
Set-Txt 20 21 "Härjedalen"
Set-Txt 20 22 "Härjedalen"
Set-Txt 20 23 "Sveg"
Set-Txt 20 25 "2023-09-07"
Set-Txt 20 26 "00:00"
Set-Txt 20 27 "2023-09-07"
Set-Txt 20 28 "00:00"
Set-Bool 20 30 $false
Set-Bool 20 31 $false
Set-Bool 20 33 $false
Set-Txt 20 49 "lennart karlsson"
Set-Txt 20 50 "lennart karlsson"
